# Update the "Metadata" sheet (sheet1): Version, Date, Publisher/Contact -> Publisher/Jurisdiction rows
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# The sheet currently has two duplicate "Contact" / "No display for ContactDetail" rows
# (rows 10 and 11), with an empty Publisher value in row 9. Remove the duplicate row,
# then turn the remaining pair into "Publisher: Alvearie Team" / "Jurisdiction: United
# States of America".
$meta.Rows.Item(11).Delete()

$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Update the "Elements" sheet (sheet2): the root Extension row's Short/Definition
# columns (K/L) now mirror the extension's own Title/Description.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Routing Number Code"
$elements.Range("L2").Value = "Customer-specific code for the ITS routing number"
